# Scheduled runner update: refresh market-price-derived columns (H:N)
# for affected Leve rows across multiple job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 664.62067
$ws.Range("I38").Value = 321.91306
$ws.Range("J38").Value = 1978.3334
$ws.Range("K38").Value = 965.7391799999999
$ws.Range("L38").Value = 5935.0002
$ws.Range("M38").Value = -593.7391799999999
$ws.Range("N38").Value = -6679.0002
$ws.Range("H39").Value = 679.9286
$ws.Range("I39").Value = 51.5
$ws.Range("J39").Value = 1151.25
$ws.Range("K39").Value = 154.5
$ws.Range("L39").Value = 3453.75
$ws.Range("M39").Value = 141.5
$ws.Range("N39").Value = -4045.75
$ws.Range("H111").Value = 998.3333
$ws.Range("I111").Value = 989.0909
$ws.Range("K111").Value = 2967.2727
$ws.Range("M111").Value = 99.72730000000001
$ws.Range("H125").Value = 2393.8
$ws.Range("I125").Value = 2191.375
$ws.Range("J125").Value = 2625.1428
$ws.Range("K125").Value = 19722.375
$ws.Range("L125").Value = 23626.2852
$ws.Range("M125").Value = -17262.375
$ws.Range("N125").Value = -28546.2852
$ws.Range("H138").Value = 3658.087
$ws.Range("I138").Value = 1958.7174
$ws.Range("J138").Value = 5357.4565
$ws.Range("K138").Value = 5876.1522
$ws.Range("L138").Value = 16072.3695
$ws.Range("M138").Value = -736.1522000000004
$ws.Range("N138").Value = -26352.3695

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14771.697
$ws.Range("I32").Value = 13102.887
$ws.Range("K32").Value = 13102.887
$ws.Range("M32").Value = -12815.887
$ws.Range("H61").Value = 4450.875
$ws.Range("I61").Value = 2183.3333
$ws.Range("J61").Value = 11253.5
$ws.Range("K61").Value = 2183.3333
$ws.Range("L61").Value = 11253.5
$ws.Range("M61").Value = -1971.3333
$ws.Range("N61").Value = -11677.5
$ws.Range("H74").Value = 2046.1578
$ws.Range("I74").Value = 1617.3125
$ws.Range("K74").Value = 1617.3125
$ws.Range("M74").Value = -743.3125
$ws.Range("H77").Value = 2046.1578
$ws.Range("I77").Value = 1617.3125
$ws.Range("K77").Value = 8086.5625
$ws.Range("M77").Value = -3718.5625
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
$ws.Range("H132").Value = 2794.6316
$ws.Range("I132").Value = 2170.2222
$ws.Range("J132").Value = 4327.273
$ws.Range("K132").Value = 6510.6666
$ws.Range("L132").Value = 12981.819
$ws.Range("M132").Value = -3980.6666
$ws.Range("N132").Value = -18041.819
$ws.Range("H136").Value = 4450.875
$ws.Range("I136").Value = 2183.3333
$ws.Range("J136").Value = 11253.5
$ws.Range("K136").Value = 6549.999899999999
$ws.Range("L136").Value = 33760.5
$ws.Range("M136").Value = -3999.999899999999
$ws.Range("N136").Value = -38860.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 29400
$ws.Range("J92").Value = 29400
$ws.Range("L92").Value = 29400
$ws.Range("N92").Value = -34392

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2692.75
$ws.Range("I105").Value = 2613.8823
$ws.Range("J105").Value = 2884.2856
$ws.Range("K105").Value = 2613.8823
$ws.Range("L105").Value = 2884.2856
$ws.Range("M105").Value = -866.8823000000002
$ws.Range("N105").Value = -6378.2856
$ws.Range("H107").Value = 1337.0938
$ws.Range("I107").Value = 1306.5238
$ws.Range("J107").Value = 1395.4546
$ws.Range("K107").Value = 1306.5238
$ws.Range("L107").Value = 1395.4546
$ws.Range("M107").Value = 613.4762000000001
$ws.Range("N107").Value = -5235.4546
$ws.Range("H132").Value = 3312.348
$ws.Range("I132").Value = 2886.625
$ws.Range("J132").Value = 4285.4287
$ws.Range("K132").Value = 8659.875
$ws.Range("L132").Value = 12856.2861
$ws.Range("M132").Value = -6129.875
$ws.Range("N132").Value = -17916.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1249.2142
$ws.Range("J113").Value = 1249.2142
$ws.Range("L113").Value = 3747.6426
$ws.Range("N113").Value = -8087.642599999999
$ws.Range("H132").Value = 2493.6333
$ws.Range("I132").Value = 1442.3334
$ws.Range("J132").Value = 3194.5
$ws.Range("K132").Value = 12981.0006
$ws.Range("L132").Value = 28750.5
$ws.Range("M132").Value = -10451.0006
$ws.Range("N132").Value = -33810.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4811.6885
$ws.Range("I122").Value = 3690.7715
$ws.Range("J122").Value = 6320.615
$ws.Range("K122").Value = 11072.3145
$ws.Range("L122").Value = 18961.845
$ws.Range("M122").Value = -8622.3145
$ws.Range("N122").Value = -23861.845
$ws.Range("H132").Value = 5422.6665
$ws.Range("I132").Value = 7784
$ws.Range("J132").Value = 3799.25
$ws.Range("K132").Value = 23352
$ws.Range("L132").Value = 11397.75
$ws.Range("M132").Value = -20822
$ws.Range("N132").Value = -16457.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3520
$ws.Range("I7").Value = 1700
$ws.Range("J7").Value = 3975
$ws.Range("K7").Value = 1700
$ws.Range("L7").Value = 3975
$ws.Range("M7").Value = -1588
$ws.Range("N7").Value = -4199
$ws.Range("H82").Value = 3212.5386
$ws.Range("I82").Value = 1690
$ws.Range("J82").Value = 3889.2222
$ws.Range("K82").Value = 1690
$ws.Range("L82").Value = 3889.2222
$ws.Range("M82").Value = -1329
$ws.Range("N82").Value = -4611.2222
$ws.Range("H85").Value = 3212.5386
$ws.Range("I85").Value = 1690
$ws.Range("J85").Value = 3889.2222
$ws.Range("K85").Value = 1690
$ws.Range("L85").Value = 3889.2222
$ws.Range("M85").Value = -442
$ws.Range("N85").Value = -6385.2222
$ws.Range("H126").Value = 3520
$ws.Range("I126").Value = 1700
$ws.Range("J126").Value = 3975
$ws.Range("K126").Value = 5100
$ws.Range("L126").Value = 11925
$ws.Range("M126").Value = -2630
$ws.Range("N126").Value = -16865
$ws.Range("H132").Value = 2730.7837
$ws.Range("I132").Value = 2132.1304
$ws.Range("J132").Value = 3714.2856
$ws.Range("K132").Value = 6396.3912
$ws.Range("L132").Value = 11142.8568
$ws.Range("M132").Value = -3866.3912
$ws.Range("N132").Value = -16202.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 49388.773
$ws.Range("I126").Value = 79743.53999999999
$ws.Range("K126").Value = 239230.62
$ws.Range("M126").Value = -236760.62
$ws.Range("H129").Value = 30001
$ws.Range("J129").Value = 30001
$ws.Range("L129").Value = 30001
$ws.Range("N129").Value = -40001
$ws.Range("H136").Value = 3098.8
$ws.Range("I136").Value = 2556.5264
$ws.Range("K136").Value = 7669.5792
$ws.Range("M136").Value = -5119.5792
